$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to Text so numeric-looking strings
# (e.g. "27.697.48", "4.960") are preserved exactly as typed, then
# drop the temporary format so cells keep their original (default) style.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.697.48'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.775.86'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '326.47'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.4627'
$ws.Range('E7').Value = '  +3.67%  '
$ws.Range('D8').Value = '0.3587'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').Value = '0.07472'
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('D10').Value = '41.98'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').Value = '1.101'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '20.81'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '6.032'
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('D15').Value = '7.243'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D16').Value = '1.774.46'
$ws.Range('E16').Value = '  +1.35%  '
$ws.Range('D17').Value = '93.58'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '0.06411'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').Value = '17.08'
$ws.Range('E21').Value = '  +1.95%  '
$ws.Range('D22').Value = '5.783'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').Value = '27.778.83'
$ws.Range('E23').Value = '  +0.91%  '
$ws.Range('D24').Value = '11.27'
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').Value = '2.079'
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('D26').Value = '164.30'
$ws.Range('E26').Value = '  +1.45%  '
$ws.Range('D27').Value = '20.35'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '1.979.56'
$ws.Range('E28').Value = '  +1.49%  '
$ws.Range('D29').Value = '2.165'
$ws.Range('D30').Value = '126.16'
$ws.Range('E30').Value = '  +1.47%  '
$ws.Range('D31').Value = '1.095'
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('D32').Value = '0.09224'
$ws.Range('E32').Value = '  +2.34%  '
$ws.Range('D33').Value = '3.675'
$ws.Range('E33').Value = '  +0.54%  '
$ws.Range('D34').Value = '5.528'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('D35').Value = '11.79'
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').Value = '0.06119'
$ws.Range('E37').Value = '  +2.60%  '
$ws.Range('D38').Value = '0.2087'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').Value = '0.6306'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').Value = '4.960'
$ws.Range('E40').Value = '  +1.04%  '
$ws.Range('D41').Value = '1.181'
$ws.Range('E41').Value = '  -1.85%  '
$ws.Range('D42').Value = '1.392'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('D43').Value = '7.771'
$ws.Range('D44').Value = '13.18'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '0.5886'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '122.23'
$ws.Range('E47').Value = '  +1.07%  '
$ws.Range('D48').Value = '1.948'
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('D49').Value = '0.06933'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('D50').Value = '1.137'
$ws.Range('E50').Value = '  -0.89%  '
$ws.Range('D51').Value = '72.27'
$ws.Range('E51').Value = '  +0.88%  '

$fmtRange.ClearFormats()
